$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 120 ("「ジュードと新しい自転車」" post) and shift the remaining
# rows (121..313) up by one, matching the target OOXML diff which removes that
# row and renumbers everything below it.
$ws.Rows.Item(120).Delete()
